$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A new year (2020) is being added to the table as column Q, mirroring the
# formatting of the previous year's column (P).
$ws.Range("P4:P5").Copy($ws.Range("Q4:Q5"))

$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 3.3

# Reflect the cell the author ended up clicking on after the edit.
$ws.Range("R4").Select()
